$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the two new cells on existing row 3 (X3, Y3)
$ws.Range("X3").Value = 0.19000099999999875
$ws.Range("Y3").Value = "Up"

# Add new row 4 with data
$ws.Range("A4").Value = 42641.890729166669
$ws.Range("A4").NumberFormat = "m/d/yy h:mm"
$ws.Range("B4").Value = 8
$ws.Range("C4").Value = "Buy"
$ws.Range("D4").Value = 18
$ws.Range("E4").Value = 9980
$ws.Range("F4").Value = 974
$ws.Range("G4").Value = 62
$ws.Range("H4").Value = 37
$ws.Range("I4").Value = 99
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 6897
$ws.Range("L4").Value = 178
$ws.Range("M4").Value = 107
$ws.Range("N4").Value = 3
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = "Bag"
$ws.Range("Q4").Value = 65.63785237683328
$ws.Range("R4").Value = 0.48
$ws.Range("S4").Value = 0.0685
$ws.Range("S4").NumberFormat = "0.00%"
$ws.Range("T4").Value = -0.0619
$ws.Range("T4").NumberFormat = "0.00%"
$ws.Range("U4").Value = 2.27
$ws.Range("V4").Value = "N/A"
$ws.Range("W4").Value = 0

$ws.Columns.AutoFit()
